$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.234.99"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.604.17"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB (numeric-looking text, force text with quote prefix)
$ws.Range("D5").Value = "'212.71"

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.08%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.38%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.33%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'18.33"
$ws.Range("E10").Value = "  +1.78%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  -0.60%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.828.13"
$ws.Range("E12").Value = "  +0.11%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.616.39"
$ws.Range("E13").Value = "  +0.90%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'4.01"
$ws.Range("E14").Value = "  +0.37%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.33%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.212.23"
$ws.Range("E16").Value = "  +0.29%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'62.02"
$ws.Range("E17").Value = "  +2.61%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +0.76%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.07%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'200.45"
$ws.Range("E20").Value = "  -1.88%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'4.25"
$ws.Range("E21").Value = "  +0.45%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -0.06%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +0.17%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +2.36%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'143.97"
$ws.Range("E25").Value = "  +1.54%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.05%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -2.22%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'15.18"
$ws.Range("E28").Value = "  -0.01%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +2.08%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +3.90%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +0.64%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.47%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -1.42%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.83%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.164.79"
$ws.Range("E36").Value = "  +4.30%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.0170"
$ws.Range("E37").Value = "  +3.73%  "

# Row 38 - PaxDollar
$ws.Range("E38").Value = "  -0.09%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "'2.32"
$ws.Range("E39").Value = "  +0.32%  "

# Row 40 - ARBITRUM
$ws.Range("D40").Value = "'0.782"
$ws.Range("E40").Value = "  +0.10%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  +0.89%  "

# Row 42 - now TrustWalletToken (was FraxShare)
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.785"
$ws.Range("E42").Value = "  +0.45%  "

# Row 43 - now FraxShare (was TrustWalletToken)
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.34"
$ws.Range("E43").Value = "  +4.22%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.739.50"
$ws.Range("E44").Value = "  +0.05%  "

# Row 45 - Quant
$ws.Range("E45").Value = "  -1.18%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +15.57%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +1.27%  "

# Row 48 - Aave
$ws.Range("D48").Value = "'54.06"
$ws.Range("E48").Value = "  +1.18%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.03%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -0.51%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  -0.16%  "
